$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the "ratings" side table (columns H:J) -------------------------
# The order cell values are assigned below matters: it determines the
# order new entries are appended to the shared-strings table, so the
# cells are touched in the exact sequence needed to reproduce that order.

# 1) New component label that seeds the shared string table addition.
$ws.Range("H5").Value2 = "Capacitor"

# 2) Column headers for the new rating columns.
$ws.Range("I1").Value2 = "Current Rating"
$ws.Range("J1").Value2 = "Voltage Rating"

# 3) Rectifier rating row.
$ws.Range("I4").Value2 = "35A"
$ws.Range("J4").Value2 = "1000V"

# 4) IGBT rating row.
$ws.Range("I2").Value2 = "30A"
$ws.Range("J2").Value2 = "600V"

# 5) Diode rating row (reuses the 30A/600V strings already created above).
$ws.Range("I3").Value2 = "30A"
$ws.Range("J3").Value2 = "600V"

# 6) Capacitor voltage rating (no current rating given).
$ws.Range("J5").Value2 = "400V"

# 7) Fill in the remaining, already-existing component labels for H2:H4.
$ws.Range("H2").Value2 = "IGBT"
$ws.Range("H3").Value2 = "Diode"
$ws.Range("H4").Value2 = "Rectifier"

# --- Add the new Heat Sink row (row 9) -----------------------------------
$ws.Range("A9").Value2 = "Heat Sink "
$ws.Range("F9").Value2 = "https://ozdisan.com/elektromekanik-komponentler/sogutucular/aluminyum-sogutucular/530802B05100G"
$ws.Range("B9").Value2 = "530802B05100G"
$ws.Range("D9").Value2 = "31,6TL"
$ws.Range("C9").Value2 = "x1"

# --- Column widths for the two new columns -------------------------------
# The underlying runtime snaps column widths to whole pixels (using a 6px
# "maximum digit width"), so the exact source widths of 21.28515625 and
# 16.7109375 characters are not bit-for-bit reproducible; these inputs
# land on the closest achievable stored widths.
$ws.Columns.Item(9).ColumnWidth = 20.5
$ws.Columns.Item(10).ColumnWidth = 15.8

# --- Selection matches the post-edit state --------------------------------
$null = $ws.Range("F18").Select()
